$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1465.462811187185
$ws.Range("C2").Value = 2458.815200492646
$ws.Range("D2").Value = 3973.532412259082
